$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 90, shifting existing rows 90-138 down to 92-140
$xlShiftDown = -4121
$ws.Rows.Item(90).Resize(2).Insert($xlShiftDown)

# Build the two new rows of data (2 rows x 18 columns, A:R)
$data = New-Object 'object[,]' 2,18

# Row 90 (index 0)
$data[0,0]  = 10
$data[0,1]  = "Vega Modelo de Temuco"
$data[0,2]  = "La Araucanía"
$data[0,3]  = 44452
$data[0,4]  = 9
$data[0,5]  = 100112039
$data[0,6]  = "Ciboulette"
$data[0,7]  = "Sin especificar"
$data[0,8]  = "Primera"
$data[0,9]  = 30
$data[0,10] = 8000
$data[0,11] = 8000
$data[0,12] = 8000
$data[0,13] = "`$/docena de atados"
$data[0,14] = "Provincia de Cautín"
$data[0,15] = 2667
$data[0,16] = 3
$data[0,17] = "Hortaliza"

# Row 91 (index 1)
$data[1,0]  = 10
$data[1,1]  = "Vega Modelo de Temuco"
$data[1,2]  = "La Araucanía"
$data[1,3]  = 44452
$data[1,4]  = 9
$data[1,5]  = 100112039
$data[1,6]  = "Ciboulette"
$data[1,7]  = "Sin especificar"
$data[1,8]  = "Segunda"
$data[1,9]  = 20
$data[1,10] = 5000
$data[1,11] = 5000
$data[1,12] = 5000
$data[1,13] = "`$/docena de atados"
$data[1,14] = "Región Metropolitana"
$data[1,15] = 1667
$data[1,16] = 3
$data[1,17] = "Hortaliza"

$ws.Range("A90:R91").Value = $data

Write-Host ("Dimension=" + $ws.UsedRange.Address())
Write-Host ("D90=" + $ws.Range("D90").Value2)
Write-Host ("D91=" + $ws.Range("D91").Value2)
Write-Host ("D92=" + $ws.Range("D92").Value2)
Write-Host ("D140=" + $ws.Range("D140").Value2)
